$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fabric Nodes")

# Remove values from POD ID column (C2:C5) while keeping formatting
$ws.Range("C2:C5").ClearContents()

# Select B2 on the active sheet (matches the <selection activeCell="B2" .../> change)
$ws.Activate()
$ws.Range("B2").Select()
